$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# New text entries for the log/dump feature (both GUI and HAL):
#   Row 11: SingleUseId11 -> Label, Center aligned, LTR, GB text "Log Data"
#   Row 12: SingleUseId12 -> Label, Center aligned, LTR, GB text "Dump Log"
$ws.Range("B11").Value = "SingleUseId11"
$ws.Range("B12").Value = "SingleUseId12"
$ws.Range("C11").Value = "Label"
$ws.Range("C12").Value = "Label"
$ws.Range("D11").Value = "Center"
$ws.Range("D12").Value = "Center"
$ws.Range("E12").Value = "Dump Log"
$ws.Range("E11").Value = "Log Data"
$ws.Range("F11").Value = "LTR"
$ws.Range("F12").Value = "LTR"
